$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldUuid = "31e1d9be-4621-4408-8307-24f412392e44"
$newUuid = "f72d2039-1438-4515-aa58-11f1078e2401"

# Update the uuid on the two existing data rows.
$ws.Range("G2").Value = $newUuid
$ws.Range("G3").Value = $newUuid

# Row 4 - duplicate of row 2's line data, with the new uuid.
$ws.Range("A4").Value = "Line:8 Stage:1"

# "01/09/2024" reads as a date to Excel's literal-assignment auto-detection.
# Copy the existing plain-text cell instead of assigning the string directly
# so the new cell stays a plain (unstyled) string, just like B2/B3.
$ws.Range("B2").Copy($ws.Range("B4"))

$ws.Range("C4").Value = "pri cl LA"

$ws.Range("D4").Value = 45300.42056299769
$ws.Range("D4").NumberFormat = $ws.Range("D2").NumberFormat

$ws.Range("E4").Value = 45300.42067873842
$ws.Range("E4").NumberFormat = $ws.Range("E2").NumberFormat

$ws.Range("F4").Value = 0.17
$ws.Range("G4").Value = $newUuid

$ws.Range("H4").Value = 10
$ws.Range("H4").NumberFormat = $ws.Range("H2").NumberFormat

# Row 5 - duplicate of row 3's line data, with the new uuid.
$ws.Range("A5").Value = "Line:8 Stage:1"

$ws.Range("B3").Copy($ws.Range("B5"))

$ws.Range("C5").Value = "pri pH HA"

$ws.Range("D5").Value = 45300.65476473379
$ws.Range("D5").NumberFormat = $ws.Range("D3").NumberFormat

$ws.Range("E5").Value = 45300.65488047454
$ws.Range("E5").NumberFormat = $ws.Range("E3").NumberFormat

$ws.Range("F5").Value = 0.17
$ws.Range("G5").Value = $newUuid

$ws.Range("H5").Value = 10
$ws.Range("H5").NumberFormat = $ws.Range("H3").NumberFormat
